# Auto update Excel log
# Appends new sensor-log rows to several sheets of the SeniorConnect
# master log workbook (mmWave, PIR, Humidity, Camera, Proximity).
#
# Values such as "2026-01-30" or "87.7%" look like dates/percentages to
# Excel's input parser, so they are written through a small helper that
# forces a Text number format before the assignment (and resets the
# style back to Normal afterwards) to guarantee the values land in the
# sheet as plain text, matching the rest of the log.

function Set-LogCell($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

function Add-LogRow($ws, $row, $date, $timestamp, $hour, $location, $value, $status) {
    Set-LogCell $ws $row 1 $date
    Set-LogCell $ws $row 2 $timestamp
    Set-LogCell $ws $row 3 $hour
    Set-LogCell $ws $row 4 $location
    Set-LogCell $ws $row 5 $value
    Set-LogCell $ws $row 6 $status
}

$wb = $excel.ActiveWorkbook

# --- mmWave: append rows 5-7 (Living Room presence detections) ---
$ws = $wb.Worksheets.Item("mmWave")
Add-LogRow $ws 5 "2026-01-30" "15:55:23" "15:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $ws 6 "2026-01-30" "15:55:34" "15:00" "Living Room" "PRESENCE_DETECTED" "Active"
Add-LogRow $ws 7 "2026-01-30" "15:55:44" "15:00" "Living Room" "PRESENCE_DETECTED" "Active"

# --- PIR: append rows 9-18 (Bathroom no-motion readings) ---
$ws = $wb.Worksheets.Item("PIR")
Add-LogRow $ws 9  "2026-01-30" "15:55:17" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 10 "2026-01-30" "15:55:19" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 11 "2026-01-30" "15:55:19" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 12 "2026-01-30" "15:55:22" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 13 "2026-01-30" "15:55:27" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 14 "2026-01-30" "15:55:32" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 15 "2026-01-30" "15:55:37" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 16 "2026-01-30" "15:55:42" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 17 "2026-01-30" "15:55:47" "15:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $ws 18 "2026-01-30" "15:55:52" "15:00" "Bathroom" "No Motion" "Inactive"

# --- Humidity: append rows 8-16 (Bathroom humidity %) ---
$ws = $wb.Worksheets.Item("Humidity")
Add-LogRow $ws 8  "2026-01-30" "15:55:17" "15:00" "Bathroom" "87.7%" "Active"
Add-LogRow $ws 9  "2026-01-30" "15:55:19" "15:00" "Bathroom" "87.7%" "Active"
Add-LogRow $ws 10 "2026-01-30" "15:55:22" "15:00" "Bathroom" "87.8%" "Active"
Add-LogRow $ws 11 "2026-01-30" "15:55:27" "15:00" "Bathroom" "87.8%" "Active"
Add-LogRow $ws 12 "2026-01-30" "15:55:32" "15:00" "Bathroom" "86.8%" "Active"
Add-LogRow $ws 13 "2026-01-30" "15:55:37" "15:00" "Bathroom" "87.8%" "Active"
Add-LogRow $ws 14 "2026-01-30" "15:55:43" "15:00" "Bathroom" "87.8%" "Active"
Add-LogRow $ws 15 "2026-01-30" "15:55:48" "15:00" "Bathroom" "87.8%" "Active"
Add-LogRow $ws 16 "2026-01-30" "15:55:53" "15:00" "Bathroom" "86.9%" "Active"

# --- Camera: append rows 4-6 (Living Room Main Door images) ---
$ws = $wb.Worksheets.Item("Camera")
Add-LogRow $ws 4 "2026-01-30" "15:55:19" "15:00" "Living Room Main Door" "Image Captured (EXIT)"  "Active"
Add-LogRow $ws 5 "2026-01-30" "15:55:24" "15:00" "Living Room Main Door" "Image Captured (ENTER)" "Active"
Add-LogRow $ws 6 "2026-01-30" "15:55:38" "15:00" "Living Room Main Door" "Image Captured (EXIT)"  "Active"

# --- Proximity: append rows 4-6 (Living Room Main Door enter/exit events) ---
$ws = $wb.Worksheets.Item("Proximity")
Add-LogRow $ws 4 "2026-01-30" "15:55:19" "15:00" "Living Room Main Door" "EXIT"  "User EXITED Living Room Main Door"
Add-LogRow $ws 5 "2026-01-30" "15:55:25" "15:00" "Living Room Main Door" "ENTER" "User ENTERED Living Room Main Door"
Add-LogRow $ws 6 "2026-01-30" "15:55:38" "15:00" "Living Room Main Door" "EXIT"  "User EXITED Living Room Main Door"
